# Login + linked-list module test data
# - Adds a "Sheet2" (invalid-login test data) after the existing "Sheet1"
# - Updates Sheet1's selection to A2:B2
# - Populates Sheet2 with the invalid-login fixture rows and formats the
#   header cells to match Sheet1's header formatting (username/password)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new worksheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Sheet2 data (column order matches the diff: A=username, B=password, C=isdatavalid) ---
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"

$ws2.Range("A2").Value = "invaliduser1"
$ws2.Range("B2").Value = "invalidpass"

$ws2.Range("B3").Value = "invalidpass2"

$ws2.Range("A4").Value = "invaliduser2"

$ws2.Range("C1").Value = "isdatavalid"
$ws2.Range("C2").Value = "N"
$ws2.Range("C3").Value = "N"
$ws2.Range("C4").Value = "N"

# --- Match the header formatting used on Sheet1 (A1/B1) ---
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws1.Range("B1").Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row heights that differ between the two sheets ---
$ws1.Rows.Item(1).RowHeight = 20.45
$ws2.Rows.Item(1).RowHeight = 16.5

# --- Selections ---
$null = $ws1.Range("A2:B2").Select()
$null = $ws2.Range("C2:C4").Select()
